# ------------------------------------------------------------------
# Template placeholder clean-up: collapse the split
# "{" / "proofErr" / "name" / "proofErr" / "}" run groups (an artifact
# of Word's auto-spellcheck breaking merge-field placeholders into
# several runs) back into a single run per placeholder, and bump the
# maintenance-fee amounts.
# ------------------------------------------------------------------

$d = $word.ActiveDocument

function Merge-Placeholder($findText) {
    $r = $d.Content
    $ok = $r.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $findText, 2)
    return $ok
}

# Simple placeholders that collapse to identical visible text (the
# Find/Replace below re-types the same characters as a single run,
# which removes the intervening <w:proofErr/> spell-check markers and
# merges the three runs Word had split the placeholder into).
Merge-Placeholder("{ownerName}") | Out-Null
Merge-Placeholder("{@ownerAddress}") | Out-Null
Merge-Placeholder("Trademark name:      {tradeMark}") | Out-Null
Merge-Placeholder("Registration Number:    {regNumber}") | Out-Null
Merge-Placeholder("Number of classes:         {classCount}") | Out-Null
Merge-Placeholder("{renewalDate}") | Out-Null
Merge-Placeholder("{markType}") | Out-Null
Merge-Placeholder("{%logoPath}") | Out-Null
Merge-Placeholder("{filingDate}") | Out-Null
Merge-Placeholder("{dateInLocation}") | Out-Null
Merge-Placeholder("{regDate}") | Out-Null
Merge-Placeholder("{intClasses}") | Out-Null
Merge-Placeholder("{serialNumber}") | Out-Null

# ------------------------------------------------------------------
# Maintenance-fee paragraph: $925 / $625 -> $1650 / $850, broken out
# into their own runs (as if the user retyped each amount in place),
# with the "_GoBack" bookmark ending up right after the new "$850".
# ------------------------------------------------------------------

# 1) Update the visible text in one pass (keeps everything in a single
#    run for now; the old bookmark collapses to a point since the
#    replaced range spans across it).
$feeRange = $d.Content
$feeRange.Find.Execute("The maintenance fee is `$925 for one class and `$625", $false, $false, $false, $false, $false, $true, 1, $false, "The maintenance fee is `$1650 for one class and `$850", 2) | Out-Null

# 2) Re-create the "_GoBack" bookmark right after the new "$850".
$afterAmount = $d.Content
$afterAmount.Find.Execute("`$850", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmPoint = $d.Range($afterAmount.End, $afterAmount.End)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

# 3) Split "The maintenance fee is $1650 for one class and $850" into
#    the target run boundaries by nudging formatting off/back on at
#    each split point (forces a new run without changing the final
#    look, matching the checked-in structure).
$prefixRange = $d.Content
$prefixRange.Find.Execute("The maintenance fee is `$", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterFirstDollar = $prefixRange.End

$num1Range = $d.Range($afterFirstDollar, $afterFirstDollar + 4)
$num1Range.Font.Bold = $false
$num1Range.Font.Bold = $true

$midRange = $d.Content
$midRange.Find.Execute(" for one class and `$", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$midRange2 = $d.Range($midRange.Start, $midRange.End)
$midRange2.Font.Bold = $false
$midRange2.Font.Bold = $true

$afterSecondDollar = $midRange.End
$num2Range = $d.Range($afterSecondDollar, $afterSecondDollar + 3)
$num2Range.Font.Bold = $false
$num2Range.Font.Bold = $true

$tailRange = $d.Content
$tailRange.Find.Execute(" for each additional class for the whole period of ten (10) years.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tailRange2 = $d.Range($tailRange.Start, $tailRange.End)
$tailRange2.Font.Bold = $false
$tailRange2.Font.Bold = $true
